$d = $word.ActiveDocument

# Locate the heading paragraph for the "Listen:" section, then the list item
# paragraph immediately following it (the one that currently just says "test").
$headingIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $headingText = $d.Paragraphs.Item($i).Range.Text
    if ($headingText -match "Listen:") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq $null) {
    throw "Could not find the 'Listen:' heading paragraph"
}

$target = $d.Paragraphs.Item($headingIndex + 1)

# Range covering the whole paragraph, including its paragraph mark, so that
# InsertXML replaces the paragraph's run content while keeping its pPr
# (style + numbering) intact.
$r = $target.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:pPr>' +
           '<w:pStyle w:val="Compact"/>' +
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1004"/></w:numPr>' +
         '</w:pPr>' +
         '<w:r><w:t xml:space="preserve">Ezra Klein Show,</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">&#8220;The Government Knows AGI is Coming&#8221;</w:t></w:r>' +
       '</w:p>'

$r.InsertXML($xml) | Out-Null
